# Update odds values that changed for the week of 2024-11-14 (FlashScore refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.63
$ws.Range("I2").Value = 2.88
$ws.Range("J2").Value = 3.5
$ws.Range("L2").Value = 3.75
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("Z2").Value = 26
$ws.Range("AA2").Value = 26
$ws.Range("AH2").Value = 7
$ws.Range("AI2").Value = 12
$ws.Range("AK2").Value = 29
$ws.Range("AN2").Value = 4.5
$ws.Range("AO2").Value = 17
$ws.Range("AP2").Value = 34
$ws.Range("AR2").Value = 101
$ws.Range("AX2").Value = 17
$ws.Range("AZ2").Value = 51
$ws.Range("BB2").Value = 301

# Row 3
$ws.Range("G3").Value = 2.2
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 3.7
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 1.91
$ws.Range("L3").Value = 4.33
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("O3").Value = 1.5
$ws.Range("P3").Value = 2.5
$ws.Range("Q3").Value = 2.5
$ws.Range("R3").Value = 1.5
$ws.Range("S3").Value = 1.57
$ws.Range("T3").Value = 2.25
$ws.Range("U3").Value = 2.1
$ws.Range("V3").Value = 1.67
$ws.Range("W3").Value = 6
$ws.Range("X3").Value = 9
$ws.Range("Y3").Value = 10
$ws.Range("Z3").Value = 21
$ws.Range("AA3").Value = 21
$ws.Range("AE3").Value = 17
$ws.Range("AH3").Value = 8.5
$ws.Range("AI3").Value = 17
$ws.Range("AJ3").Value = 15
$ws.Range("AK3").Value = 41
$ws.Range("AL3").Value = 34
$ws.Range("AM3").Value = 41
$ws.Range("AN3").Value = 4
$ws.Range("AO3").Value = 13
$ws.Range("AP3").Value = 29
$ws.Range("AR3").Value = 81
$ws.Range("AT3").Value = 2.25
$ws.Range("AU3").Value = 9
$ws.Range("AW3").Value = 5.5
$ws.Range("AX3").Value = 21
$ws.Range("AY3").Value = 34
$ws.Range("AZ3").Value = 81
$ws.Range("BA3").Value = 126

# Row 4
$ws.Range("G4").Value = 4.1
$ws.Range("H4").Value = 2.9
$ws.Range("I4").Value = 2.1
$ws.Range("J4").Value = 4.75
$ws.Range("L4").Value = 2.88
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("W4").Value = 9
$ws.Range("X4").Value = 19
$ws.Range("AD4").Value = 6
$ws.Range("AH4").Value = 5.5
$ws.Range("AI4").Value = 8.5
$ws.Range("AK4").Value = 19
$ws.Range("AX4").Value = 12
$ws.Range("BB4").Value = 251

# Row 5
$ws.Range("G5").Value = 1.3
$ws.Range("H5").Value = 4.75
$ws.Range("I5").Value = 12
$ws.Range("J5").Value = 1.8
$ws.Range("L5").Value = 10
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 9
$ws.Range("Q5").Value = 2
$ws.Range("R5").Value = 1.85
$ws.Range("AC5").Value = 9
$ws.Range("AD5").Value = 10
$ws.Range("AH5").Value = 21
$ws.Range("AI5").Value = 51
$ws.Range("AJ5").Value = 34
$ws.Range("AO5").Value = 6
$ws.Range("AQ5").Value = 17

# Row 6
$ws.Range("G6").Value = 1.48
$ws.Range("H6").Value = 4.1
$ws.Range("I6").Value = 7.5
$ws.Range("J6").Value = 2.05
$ws.Range("L6").Value = 7
$ws.Range("U6").Value = 2.2
$ws.Range("V6").Value = 1.62
$ws.Range("X6").Value = 6
$ws.Range("Z6").Value = 9.5
$ws.Range("AB6").Value = 34
$ws.Range("AH6").Value = 17
$ws.Range("AI6").Value = 41
$ws.Range("AJ6").Value = 23
$ws.Range("AK6").Value = 81
$ws.Range("AN6").Value = 3.25
$ws.Range("AS6").Value = 201
$ws.Range("AU6").Value = 10
$ws.Range("AW6").Value = 8
$ws.Range("AX6").Value = 41
$ws.Range("BA6").Value = 201

# Row 8
$ws.Range("M8").Value = 1.04
$ws.Range("N8").Value = 13

# Row 9
$ws.Range("G9").Value = 1.3
$ws.Range("H9").Value = 4.75
$ws.Range("I9").Value = 11
$ws.Range("J9").Value = 1.83
$ws.Range("K9").Value = 2.38
$ws.Range("L9").Value = 9
$ws.Range("M9").Value = 1.05
$ws.Range("N9").Value = 11
$ws.Range("O9").Value = 1.29
$ws.Range("P9").Value = 3.5
$ws.Range("Q9").Value = 1.95
$ws.Range("R9").Value = 1.9
$ws.Range("U9").Value = 2.38
$ws.Range("V9").Value = 1.53
$ws.Range("X9").Value = 5.5
$ws.Range("Z9").Value = 7.5
$ws.Range("AC9").Value = 9.5
$ws.Range("AD9").Value = 9.5
$ws.Range("AE9").Value = 26
$ws.Range("AF9").Value = 101
$ws.Range("AH9").Value = 21
$ws.Range("AI9").Value = 51
$ws.Range("AJ9").Value = 34
$ws.Range("AK9").Value = 151
$ws.Range("AL9").Value = 81
$ws.Range("AM9").Value = 81
$ws.Range("AN9").Value = 3.1
$ws.Range("AO9").Value = 6.5
$ws.Range("AQ9").Value = 19
$ws.Range("AU9").Value = 11
$ws.Range("AW9").Value = 10
$ws.Range("AX9").Value = 51
$ws.Range("AY9").Value = 51
$ws.Range("AZ9").Value = 251
$ws.Range("BA9").Value = 301

# Row 11
$ws.Range("M11").Value = 1.06
$ws.Range("N11").Value = 10

# The Uruguay - Primera Division (Wanderers vs Penarol) match that was row 14 is no
# longer part of this week's fixtures, so remove the entire row (shrinks used range
# from A1:BD14 to A1:BD13).
$ws.Rows.Item(14).Delete()
